$d = $word.ActiveDocument

# 1. Retitle the second heading paragraph from "General Remarks" to "Run Model".
#    Delete the whole paragraph (runs + proofErr markers) and insert a fresh
#    one carrying the same style/formatting, so no stray <w:proofErr/> is left
#    behind from the old "General " / "Remarks" run split.
$headingPara = $d.Paragraphs(2)
$headingPara.Range.Delete()

$bookmarkPara = $d.Paragraphs(2)
$bookmarkPara.Range.InsertParagraphBefore()
$newHeading = $d.Paragraphs(2)
$newHeading.Range.Text = "Run Model"
$newHeading.Style = "Kop1"

# 2. Append a new, empty paragraph right after the bookmark paragraph
#    (before the sectPr), keeping the nl-BE language mark but no run.
$bookmarkPara = $d.Paragraphs(3)
$bookmarkPara.Range.InsertParagraphAfter()
$newEmpty = $d.Paragraphs(4)
$newEmpty.Range.Text = "x"
$r = $d.Paragraphs(4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = ""
